$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (row number -> column letter -> new value), taken from
# the recomputed NATMI LR-pair statistics for the App-Cav1 pair after
# incorporating additional replicate data (Dr Hou's advice: ligand- and
# receptor-expressing cell counts go from 1 to 3, with updated total/average
# expression values and re-derived specificity + edge-weight columns).
$data = @{
  2 = @{ E=3; G=103.4275383333333; H=310.282615; I=0.2485530285127421; J=0.2485530285127421; K=3; M=833.4308676666666; N=2500.292603; O=0.8518935545813505; P=0.8518935545813505; Q=86199.70301377743; R=775797.3271239968; S=0.2117407229616796; T=0.2117407229616796 }
  3 = @{ E=3; G=103.4275383333333; H=310.282615; I=0.2485530285127421; J=0.2485530285127421; K=3; M=17.73945766666667; N=53.218373; O=0.01813243333584592; P=0.01813243333584592; Q=1834.748437831711; R=16512.7359404854; S=0.004506871219929907; T=0.004506871219929906 }
  4 = @{ E=3; G=103.4275383333333; H=310.282615; I=0.2485530285127421; J=0.2485530285127421; K=3; M=1.344749666666667; N=4.034249; O=0.001374539410528448; P=0.001374539410528448; Q=139.0841476979039; R=1251.757329281135; S=0.0003416459332969651; T=0.000341645933296965 }
  5 = @{ E=3; G=103.4275383333333; H=310.282615; I=0.2485530285127421; J=0.2485530285127421; K=3; M=125.812397; N=377.437191; O=0.1285994726722751; P=0.1285994726722751; Q=13012.46651352605; R=117112.1986217345; S=0.0319637883978356; T=0.0319637883978356 }
  6 = @{ E=3; G=216.130539; H=648.391617; I=0.5193964865470273; J=0.5193964865470272; K=3; M=833.4308676666666; N=2500.292603; O=0.8518935545813505; P=0.8518935545813505; Q=180129.8626480343; R=1621168.763832309; S=0.4424705191616117; T=0.4424705191616116 }
  7 = @{ E=3; G=216.130539; H=648.391617; I=0.5193964865470273; J=0.5193964865470272; K=3; M=17.73945766666667; N=53.218373; O=0.01813243333584592; P=0.01813243333584592; Q=3834.038547064349; R=34506.34692357914; S=0.009417922167186565; T=0.009417922167186564 }
  8 = @{ E=3; G=216.130539; H=648.391617; I=0.5193964865470273; J=0.5193964865470272; K=3; M=1.344749666666667; N=4.034249; O=0.001374539410528448; P=0.001374539410528448; Q=290.641470276737; R=2615.773232490633; S=0.000713930940448898; T=0.0007139309404488977 }
  9 = @{ E=3; G=216.130539; H=648.391617; I=0.5193964865470273; J=0.5193964865470272; K=3; M=125.812397; N=377.437191; O=0.1285994726722751; P=0.1285994726722751; Q=27191.90117649199; R=244727.1105884279; S=0.06679411427778016; T=0.06679411427778015 }
  10 = @{ E=3; G=71.607325; H=214.821975; I=0.1720839321833696; J=0.1720839321833696; K=3; M=833.4308676666666; N=2500.292603; O=0.8518935545813505; P=0.8518935545813505; Q=59679.75500603899; R=537117.7950543509; S=0.1465971926740268; T=0.1465971926740268 }
  11 = @{ E=3; G=71.607325; H=214.821975; I=0.1720839321833696; J=0.1720839321833696; K=3; M=17.73945766666667; N=53.218373; O=0.01813243333584592; P=0.01813243333584592; Q=1270.275110460742; R=11432.47599414667; S=0.003120300428485179; T=0.003120300428485179 }
  12 = @{ E=3; G=71.607325; H=214.821975; I=0.1720839321833696; J=0.1720839321833696; K=3; M=1.344749666666667; N=4.034249; O=0.001374539410528448; P=0.001374539410528448; Q=96.29392642464167; R=866.645337821775; S=0.0002365361467047463; T=0.0002365361467047462 }
  13 = @{ E=3; G=71.607325; H=214.821975; I=0.1720839321833696; J=0.1720839321833696; K=3; M=125.812397; N=377.437191; O=0.1285994726722751; P=0.1285994726722751; Q=9009.089201008026; R=81081.80280907224; S=0.02212990293415288; T=0.02212990293415288 }
  14 = @{ E=3; G=24.953198; H=74.859594; I=0.05996655275686102; J=0.05996655275686102; K=3; M=833.4308676666666; N=2500.292603; O=0.8518935545813505; P=0.8518935545813505; Q=20796.76546019813; R=187170.8891417832; S=0.05108511978403242; T=0.05108511978403241 }
  15 = @{ E=3; G=24.953198; H=74.859594; I=0.05996655275686102; J=0.05996655275686102; K=3; M=17.73945766666667; N=53.218373; O=0.01813243333584592; P=0.01813243333584592; Q=442.6561995689514; R=3983.905796120562; S=0.00108733952024427; T=0.00108733952024427 }
  16 = @{ E=3; G=24.953198; H=74.859594; I=0.05996655275686102; J=0.05996655275686102; K=3; M=1.344749666666667; N=4.034249; O=0.001374539410528448; P=0.001374539410528448; Q=33.55580469276734; R=302.002242234906; S=0.00008242639007783884; T=0.00008242639007783882 }
  17 = @{ E=3; G=24.953198; H=74.859594; I=0.05996655275686102; J=0.05996655275686102; K=3; M=125.812397; N=377.437191; O=0.1285994726722751; P=0.1285994726722751; Q=3139.421653195607; R=28254.79487876046; S=0.007711667062506495; T=0.007711667062506494 }
}

$colMap = @{ 'E'=5; 'F'=6; 'G'=7; 'H'=8; 'I'=9; 'J'=10; 'K'=11; 'L'=12; 'M'=13; 'N'=14; 'O'=15; 'P'=16; 'Q'=17; 'R'=18; 'S'=19; 'T'=20 }

foreach ($r in $data.Keys) {
    $rowData = $data[$r]
    foreach ($col in $rowData.Keys) {
        $c = $colMap[$col]
        $ws.Cells.Item([int]$r, $c).Value = $rowData[$col]
    }
}

Write-Output "Updated $($data.Keys.Count) rows"
